$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-03"
$ws.Cells.Item($row, 2).Value = "14:18:36"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"
$ws.Cells.Item($row, 5).Value = 139859
$ws.Cells.Item($row, 6).Value = 142849
$ws.Cells.Item($row, 7).Value = 171793
$ws.Cells.Item($row, 8).Value = 146692
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 117198
$ws.Cells.Item($row, 11).Value = 223830
$ws.Cells.Item($row, 12).Value = 248042
$ws.Cells.Item($row, 13).Value = 183891
$ws.Cells.Item($row, 14).Value = 109875
$ws.Cells.Item($row, 15).Value = 39935
$ws.Cells.Item($row, 16).Value = 30847
$ws.Cells.Item($row, 17).Value = 72116
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41288
$ws.Cells.Item($row, 20).Value = -1
